$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column Q (rows 4-34, the year-data rows) into the new column R so
# that the new column inherits the exact same per-row cell styles (and,
# incidentally, correct shared-string typed cells) as column Q already
# uses -- this avoids creating brand-new (duplicate) style/font records.
$ws.Range("Q4:Q34").Copy($ws.Range("R4:R34"))

# Now overwrite the copied values in column R with the real 2021 data
# taken from the edit. Text "-" cells are written as text, everything
# else as numbers.
$ws.Range("R4").Value = 2021

$ws.Range("R5").Value = 0.8
$ws.Range("R6").Value = 0.4
$ws.Range("R7").Value = 1.2
$ws.Range("R8").Value = 0.2
$ws.Range("R9").Value = "-"
$ws.Range("R10").Value = 0.4
$ws.Range("R11").Value = 0.6
$ws.Range("R12").Value = 0.8
$ws.Range("R13").Value = 0.5
$ws.Range("R14").Value = 0.4
$ws.Range("R15").Value = "-"
$ws.Range("R16").Value = 0.8
$ws.Range("R17").Value = 0.3
$ws.Range("R18").Value = 0.7
$ws.Range("R19").Value = "-"
$ws.Range("R20").Value = 0.5
$ws.Range("R21").Value = 0.1
$ws.Range("R22").Value = 0.8
$ws.Range("R23").Value = 1.1000000000000001
$ws.Range("R24").Value = 1.5
$ws.Range("R25").Value = 0.7
$ws.Range("R26").Value = 2.2000000000000002
$ws.Range("R27").Value = 1
$ws.Range("R28").Value = 3.5
$ws.Range("R29").Value = 0.8
$ws.Range("R30").Value = 0.2
$ws.Range("R31").Value = 1.6
$ws.Range("R32").Value = 0.3
$ws.Range("R33").Value = "-"
$ws.Range("R34").Value = 0.6

# Move the active selection like the author's recorded session.
$ws.Range("Q11").Select()
